# Investigation report template: bold the three field labels (splitting
# each into its own run), and normalize "{{ placeholder }}" spacing
# throughout the document.

$d = $word.ActiveDocument

# --- Part 1: label fields that get split into a bold "label:" run plus
#     a plain "{{ var }}..." run. -----------------------------------------
$labelFields = @(
    @("贷款用途：{{loan_use}}", "贷款用途：", "{{ loan_use }}"),
    @("贷款金额：{{loan_amount}}元（{{loan_amount_cn}}）", "贷款金额：", "{{ loan_amount }}元（{{ loan_amount_cn }}）"),
    @("贷款期限：{{loan_term}}个月", "贷款期限：", "{{ loan_term }}个月")
)

foreach ($field in $labelFields) {
    $oldText = $field[0]
    $label = $field[1]
    $rest = $field[2]

    $rng = $d.Content
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)

    if ($rng.Find.Found) {
        # Rewrite the run's text (label + normalized placeholder), then
        # carve out the leading "label：" span into its own bold run by
        # applying Font.Bold to just that sub-range.
        $rng.Text = $label + $rest

        $labelRange = $d.Range($rng.Start, $rng.Start + $label.Length)
        $labelRange.Font.Bold = 1
    }
}

# --- Part 2: plain placeholder spacing normalization ----------------------
$plainFields = @(
    @("{{main_borrower_summary}}", "{{ main_borrower_summary }}"),
    @("{{joint_borrowers_summary}}", "{{ joint_borrowers_summary }}"),
    @("{{guarantors_summary}}", "{{ guarantors_summary }}"),
    @("{{collaterals_summary}}", "{{ collaterals_summary }}")
)

foreach ($field in $plainFields) {
    $rng = $d.Content
    $rng.Find.Execute($field[0], $true, $false, $false, $false, $false, `
                       $true, 1, $false, $field[1], 2)
}
